$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header row (row 1) was missing the "vehicle-L" header that every data
# row below it already uses in column C ("Vehicle", "Model_scores", ...).
# Re-add it so column C is correctly labelled, same as column B ("vehicle-H").
$ws.Range("C1").Value = "vehicle-L"

# Column A now holds noticeably long labels (e.g. the model-score field
# names); auto-fit it so the text isn't clipped, matching what Excel does
# automatically when the sheet is touched and re-saved.
$ws.Columns("A:A").AutoFit()

# Restore the selection that was active when the workbook was last saved.
$ws.Range("A6").Select()
